$wb = $excel.ActiveWorkbook

# --- optimization_parameters sheet: restructure rows ---
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 had duplicate "value" cells in C1:F1 left over from a copy/paste; remove them.
$ws.Range("C1:F1").ClearContents()

# Insert a new row after the "Model"/"Sigmoid" row (row 8) for the new
# "L_curve" parameter, then relabel "Model" -> "production_function".
$ws.Range("A8").Value = "production_function"
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").Style = $ws.Range("B2").Style

# The old "Deletion" row (originally row 16, now row 17 after the insert
# above) is no longer used; remove it entirely.
$ws.Rows.Item(17).Delete()

# This sheet is now the active / selected tab.
$ws.Activate()
$ws.Range("C1:K3").Select()
